# Revert "Updated example to have a scale column"
# The "Scale" column (D) that was previously added to Sheet1's position
# table is removed again: the header in D9 and the value 1 in D10:D64.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Remove the whole "Scale" column content (header + the 1's below it).
$ws.Range("D9:D64").ClearContents() | Out-Null

# Restore the selection that was active before the scale column existed.
$ws.Range("I28").Select() | Out-Null
